$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.139.00'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.54%  '
$ws.Range("E2").NumberFormat = "General"
$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.860.98'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.43%  '
$ws.Range("E3").NumberFormat = "General"
$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9995'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("E4").NumberFormat = "General"
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.67'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.35%  '
$ws.Range("E5").NumberFormat = "General"
$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.13%  '
$ws.Range("E6").NumberFormat = "General"
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4696'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E7").NumberFormat = "General"
$ws.Range("B8").Value = 'OKB'
$ws.Range("C8").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '42.78'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.34%  '
$ws.Range("E8").NumberFormat = "General"
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2863'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.18%  '
$ws.Range("E9").NumberFormat = "General"
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06471'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.56%  '
$ws.Range("E10").NumberFormat = "General"
$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.59'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -4.62%  '
$ws.Range("E11").NumberFormat = "General"
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07671'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -4.33%  '
$ws.Range("E12").NumberFormat = "General"
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.850.96'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.07%  '
$ws.Range("E13").NumberFormat = "General"
$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '93.70'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.30%  '
$ws.Range("E14").NumberFormat = "General"
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.059'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.07%  '
$ws.Range("E15").NumberFormat = "General"
$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.6792'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.81%  '
$ws.Range("E16").NumberFormat = "General"
$ws.Range("B17").Value = 'BitcoinCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '268.62'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.05%  '
$ws.Range("E17").NumberFormat = "General"
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '30.123.22'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.65%  '
$ws.Range("E18").NumberFormat = "General"
$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.30'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -5.18%  '
$ws.Range("E19").NumberFormat = "General"
$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9997'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("E20").NumberFormat = "General"
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000007511'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.31%  '
$ws.Range("E21").NumberFormat = "General"
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.100.75'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.68%  '
$ws.Range("E22").NumberFormat = "General"
$ws.Range("B23").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9991'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("E23").NumberFormat = "General"
$ws.Range("B24").Value = 'Uniswap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.161'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.02%  '
$ws.Range("E24").NumberFormat = "General"
$ws.Range("B25").Value = 'Chainlink'
$ws.Range("C25").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.099'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.87%  '
$ws.Range("E25").NumberFormat = "General"
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.308'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.89%  '
$ws.Range("E26").NumberFormat = "General"
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '166.43'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.72%  '
$ws.Range("E27").NumberFormat = "General"
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.72'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.85%  '
$ws.Range("E28").NumberFormat = "General"
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.878'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.68%  '
$ws.Range("E29").NumberFormat = "General"
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.371'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.24%  '
$ws.Range("E30").NumberFormat = "General"
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09814'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.65%  '
$ws.Range("E31").NumberFormat = "General"
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.498'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.45%  '
$ws.Range("E32").NumberFormat = "General"
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.210'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.40%  '
$ws.Range("E33").NumberFormat = "General"
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.992'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.77%  '
$ws.Range("E34").NumberFormat = "General"
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.04666'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.93%  '
$ws.Range("E35").NumberFormat = "General"
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.108'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.48%  '
$ws.Range("E36").NumberFormat = "General"
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6841'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.23%  '
$ws.Range("E37").NumberFormat = "General"
$ws.Range("B38").Value = 'HuobiToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.712'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.04%  '
$ws.Range("E38").NumberFormat = "General"
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01821'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.13%  '
$ws.Range("E39").NumberFormat = "General"
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.717'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.63%  '
$ws.Range("E40").NumberFormat = "General"
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.348'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.40%  '
$ws.Range("E41").NumberFormat = "General"
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '70.09'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.57%  '
$ws.Range("E42").NumberFormat = "General"
$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9992'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("E43").NumberFormat = "General"
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8326'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.07%  '
$ws.Range("E44").NumberFormat = "General"
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.880'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.89%  '
$ws.Range("E45").NumberFormat = "General"
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.86'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.86%  '
$ws.Range("E46").NumberFormat = "General"
$ws.Range("B47").Value = 'TheSandbox'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4040'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.95%  '
$ws.Range("E47").NumberFormat = "General"
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.185'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.06%  '
$ws.Range("E48").NumberFormat = "General"
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '925.24'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.40%  '
$ws.Range("E49").NumberFormat = "General"
$ws.Range("B50").Value = 'Aptos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.910'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.04%  '
$ws.Range("E50").NumberFormat = "General"
$ws.Range("B51").Value = 'Elrond'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '34.13'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.91%  '
$ws.Range("E51").NumberFormat = "General"
